$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    7 = @('aa', 'Agree/Accept')
    14 = @('aa', 'Agree/Accept')
    15 = @('sd', 'Statement-non-opinion')
    16 = @('sd', 'Statement-non-opinion')
    21 = @('aa', 'Agree/Accept')
    42 = @('sv', 'Statement-opinion')
    48 = @('ba', 'Appreciation')
    50 = @('sd', 'Statement-non-opinion')
    66 = @('sv', 'Statement-opinion')
    92 = @('sd', 'Statement-non-opinion')
    99 = @('aa', 'Agree/Accept')
    100 = @('sd', 'Statement-non-opinion')
    108 = @('aa', 'Agree/Accept')
    127 = @('sd', 'Statement-non-opinion')
    128 = @('sd', 'Statement-non-opinion')
    138 = @('b', 'Acknowledge (Backchannel)')
    144 = @('sv', 'Statement-opinion')
    154 = @('b', 'Acknowledge (Backchannel)')
    165 = @('b', 'Acknowledge (Backchannel)')
    170 = @('sd', 'Statement-non-opinion')
    174 = @('b', 'Acknowledge (Backchannel)')
    178 = @('sv', 'Statement-opinion')
    184 = @('sv', 'Statement-opinion')
    188 = @('sd', 'Statement-non-opinion')
    189 = @('sv', 'Statement-opinion')
    192 = @('aa', 'Agree/Accept')
    196 = @('b', 'Acknowledge (Backchannel)')
    199 = @('aa', 'Agree/Accept')
    216 = @('b', 'Acknowledge (Backchannel)')
    233 = @('sv', 'Statement-opinion')
    239 = @('aa', 'Agree/Accept')
    243 = @('sv', 'Statement-opinion')
    265 = @('sd', 'Statement-non-opinion')
    268 = @('aa', 'Agree/Accept')
    285 = @('sv', 'Statement-opinion')
    327 = @('aa', 'Agree/Accept')
    340 = @('sd', 'Statement-non-opinion')
    347 = @('sd', 'Statement-non-opinion')
    351 = @('ba', 'Appreciation')
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
